$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JL PCBA")

# R49 used to be lumped into the 47k group on row 24
# (R10,R25,R28,R41,R44,R49,R62). Per the commit message, R49 actually
# belongs to a 4.7k (4K7) part, so split it out of that group...
$ws.Range("B24").Value = "R10,R25,R28,R41,R44,R62"

# ...and insert a dedicated row for it right below, copying row 24's
# formatting (thin border) first so the new row isn't left unformatted.
$ws.Rows.Item(25).Insert()
$ws.Range("A24:D24").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)

$ws.Range("A25").Value = "4.7k"
$ws.Range("B25").Value = "R49"
$ws.Range("C25").Value = "R_0805_HandSoldering"
$ws.Range("D25").Value = "C23162"

# Call out the change in red, like the author did.
$ws.Range("A25:D25").Font.Color = 255

# The BOM named range covered rows 1-40; it now needs to stretch one row
# further to keep including everything through the new last row (41).
$wb.Names.Item("FreeDSP_SMD_AB_plus_BOM").RefersTo = "='JL PCBA'!`$B`$1:`$D`$41"

$ws.Range("C18").Select()
